$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 549.8461
$ws.Range("I33").Value = 163.59091
$ws.Range("K33").Value = 163.59091
$ws.Range("M33").Value = 65.40908999999999
$ws.Range("H74").Value = 3629.8096
$ws.Range("I74").Value = 3675.0527
$ws.Range("J74").Value = 3200
$ws.Range("K74").Value = 3675.0527
$ws.Range("L74").Value = 3200
$ws.Range("M74").Value = -2739.0527
$ws.Range("N74").Value = -5072
$ws.Range("H76").Value = 37077220
$ws.Range("I76").Value = 38503132
$ws.Range("J76").Value = 3500
$ws.Range("K76").Value = 38503132
$ws.Range("L76").Value = 3500
$ws.Range("M76").Value = -38502817
$ws.Range("N76").Value = -4130
$ws.Range("H77").Value = 3629.8096
$ws.Range("I77").Value = 3675.0527
$ws.Range("J77").Value = 3200
$ws.Range("K77").Value = 18375.2635
$ws.Range("L77").Value = 16000
$ws.Range("M77").Value = -13695.2635
$ws.Range("N77").Value = -25360
$ws.Range("H79").Value = 37077220
$ws.Range("I79").Value = 38503132
$ws.Range("J79").Value = 3500
$ws.Range("K79").Value = 38503132
$ws.Range("L79").Value = 3500
$ws.Range("M79").Value = -38502040
$ws.Range("N79").Value = -5684
$ws.Range("H93").Value = 27000
$ws.Range("J93").Value = 27000
$ws.Range("L93").Value = 27000
$ws.Range("N93").Value = -31992
$ws.Range("H133").Value = 17600.2
$ws.Range("J133").Value = 17600.2
$ws.Range("L133").Value = 17600.2
$ws.Range("N133").Value = -27720.2
$ws.Range("H141").Value = 2448.3333
$ws.Range("I141").Value = 2129.375
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 6388.125
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = -1208.125
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5500.86
$ws.Range("I32").Value = 5213.6807
$ws.Range("K32").Value = 5213.6807
$ws.Range("M32").Value = -4926.6807
$ws.Range("H61").Value = 3928.6365
$ws.Range("I61").Value = 3733.5
$ws.Range("J61").Value = 4162.8
$ws.Range("K61").Value = 3733.5
$ws.Range("L61").Value = 4162.8
$ws.Range("M61").Value = -3521.5
$ws.Range("N61").Value = -4586.8
$ws.Range("H63").Value = 2927.2727
$ws.Range("I63").Value = 2911.111
$ws.Range("J63").Value = 3000
$ws.Range("K63").Value = 2911.111
$ws.Range("L63").Value = 3000
$ws.Range("M63").Value = -2225.111
$ws.Range("N63").Value = -4372
$ws.Range("H66").Value = 2927.2727
$ws.Range("I66").Value = 2911.111
$ws.Range("J66").Value = 3000
$ws.Range("K66").Value = 14555.555
$ws.Range("L66").Value = 15000
$ws.Range("M66").Value = -11123.555
$ws.Range("N66").Value = -21864
$ws.Range("H88").Value = 15798.143
$ws.Range("I88").Value = 1428.5
$ws.Range("J88").Value = 26575.375
$ws.Range("K88").Value = 1428.5
$ws.Range("L88").Value = 26575.375
$ws.Range("M88").Value = -1022.5
$ws.Range("N88").Value = -27387.375
$ws.Range("H91").Value = 15798.143
$ws.Range("I91").Value = 1428.5
$ws.Range("J91").Value = 26575.375
$ws.Range("K91").Value = 1428.5
$ws.Range("L91").Value = 26575.375
$ws.Range("M91").Value = -24.5
$ws.Range("N91").Value = -29383.375
$ws.Range("H132").Value = 213311
$ws.Range("I132").Value = 33666.973
$ws.Range("J132").Value = 628737.8
$ws.Range("K132").Value = 101000.919
$ws.Range("L132").Value = 1886213.4
$ws.Range("M132").Value = -98470.91899999999
$ws.Range("N132").Value = -1891273.4
$ws.Range("H136").Value = 3928.6365
$ws.Range("I136").Value = 3733.5
$ws.Range("J136").Value = 4162.8
$ws.Range("K136").Value = 11200.5
$ws.Range("L136").Value = 12488.4
$ws.Range("M136").Value = -8650.5
$ws.Range("N136").Value = -17588.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 21074.777
$ws.Range("J81").Value = 21074.777
$ws.Range("L81").Value = 21074.777
$ws.Range("N81").Value = -23196.777
$ws.Range("H84").Value = 21074.777
$ws.Range("J84").Value = 21074.777
$ws.Range("L84").Value = 63224.33099999999
$ws.Range("N84").Value = -73832.33099999999
$ws.Range("H92").Value = 22100
$ws.Range("J92").Value = 22100
$ws.Range("L92").Value = 22100
$ws.Range("N92").Value = -27092
$ws.Range("H105").Value = 2654633.2
$ws.Range("I105").Value = 2654633.2
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2654633.2
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -2652886.2
$ws.Range("N105").ClearContents()
$ws.Range("H134").Value = 26344158
$ws.Range("I134").Value = 2054.6128
$ws.Range("J134").Value = 143002050
$ws.Range("K134").Value = 6163.8384
$ws.Range("L134").Value = 429006150
$ws.Range("M134").Value = -3628.8384
$ws.Range("N134").Value = -429011220

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26409.334
$ws.Range("I31").Value = 35978.414
$ws.Range("J31").Value = 13795.546
$ws.Range("K31").Value = 35978.414
$ws.Range("L31").Value = 13795.546
$ws.Range("M31").Value = -35683.414
$ws.Range("N31").Value = -14385.546
$ws.Range("H32").Value = 10627.5
$ws.Range("I32").Value = 1255
$ws.Range("J32").Value = 20000
$ws.Range("K32").Value = 1255
$ws.Range("L32").Value = 20000
$ws.Range("M32").Value = -939
$ws.Range("N32").Value = -20632
$ws.Range("H34").Value = 26409.334
$ws.Range("I34").Value = 35978.414
$ws.Range("J34").Value = 13795.546
$ws.Range("K34").Value = 35978.414
$ws.Range("L34").Value = 13795.546
$ws.Range("M34").Value = -35776.414
$ws.Range("N34").Value = -14199.546
$ws.Range("H58").Value = 2774.3257
$ws.Range("I58").Value = 961.3913
$ws.Range("J58").Value = 4859.2
$ws.Range("K58").Value = 961.3913
$ws.Range("L58").Value = 4859.2
$ws.Range("M58").Value = -758.3913
$ws.Range("N58").Value = -5265.2
$ws.Range("H62").Value = 2435
$ws.Range("I62").Value = 2462
$ws.Range("J62").Value = 2300
$ws.Range("K62").Value = 2462
$ws.Range("L62").Value = 2300
$ws.Range("M62").Value = -1838
$ws.Range("N62").Value = -3548
$ws.Range("H65").Value = 2435
$ws.Range("I65").Value = 2462
$ws.Range("J65").Value = 2300
$ws.Range("K65").Value = 12310
$ws.Range("L65").Value = 11500
$ws.Range("M65").Value = -9190
$ws.Range("N65").Value = -17740
$ws.Range("H96").Value = 21980
$ws.Range("J96").Value = 21980
$ws.Range("L96").Value = 21980
$ws.Range("N96").Value = -27472
$ws.Range("H105").Value = 1502.5
$ws.Range("I105").Value = 1403.3334
$ws.Range("J105").Value = 1800
$ws.Range("K105").Value = 1403.3334
$ws.Range("L105").Value = 1800
$ws.Range("M105").Value = 343.6666
$ws.Range("N105").Value = -5294
$ws.Range("H136").Value = 2774.3257
$ws.Range("I136").Value = 961.3913
$ws.Range("J136").Value = 4859.2
$ws.Range("K136").Value = 2884.1739
$ws.Range("L136").Value = 14577.6
$ws.Range("M136").Value = -334.1738999999998
$ws.Range("N136").Value = -19677.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 17547732
$ws.Range("I75").Value = 1543
$ws.Range("J75").Value = 22226716
$ws.Range("K75").Value = 4629
$ws.Range("L75").Value = 66680148
$ws.Range("M75").Value = -3631
$ws.Range("N75").Value = -66682144
$ws.Range("H78").Value = 17547732
$ws.Range("I78").Value = 1543
$ws.Range("J78").Value = 22226716
$ws.Range("K78").Value = 13887
$ws.Range("L78").Value = 200040444
$ws.Range("M78").Value = -8895
$ws.Range("N78").Value = -200050428

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4080.6
$ws.Range("I70").Value = 4000
$ws.Range("J70").Value = 4201.5
$ws.Range("K70").Value = 4000
$ws.Range("L70").Value = 4201.5
$ws.Range("M70").Value = -3730
$ws.Range("N70").Value = -4741.5
$ws.Range("H73").Value = 4080.6
$ws.Range("I73").Value = 4000
$ws.Range("J73").Value = 4201.5
$ws.Range("K73").Value = 4000
$ws.Range("L73").Value = 4201.5
$ws.Range("M73").Value = -3064
$ws.Range("N73").Value = -6073.5
$ws.Range("H80").Value = 7186.3184
$ws.Range("I80").Value = 2920.75
$ws.Range("J80").Value = 12305
$ws.Range("K80").Value = 2920.75
$ws.Range("L80").Value = 12305
$ws.Range("M80").Value = -1922.75
$ws.Range("N80").Value = -14301
$ws.Range("H83").Value = 7186.3184
$ws.Range("I83").Value = 2920.75
$ws.Range("J83").Value = 12305
$ws.Range("K83").Value = 14603.75
$ws.Range("L83").Value = 61525
$ws.Range("M83").Value = -9611.75
$ws.Range("N83").Value = -71509
$ws.Range("H132").Value = 73886.28999999999
$ws.Range("I132").Value = 2222.4443
$ws.Range("J132").Value = 202881.2
$ws.Range("K132").Value = 6667.3329
$ws.Range("L132").Value = 608643.6000000001
$ws.Range("M132").Value = -4137.3329
$ws.Range("N132").Value = -613703.6000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 12608.889
$ws.Range("I68").Value = 13812.5
$ws.Range("J68").Value = 2980
$ws.Range("K68").Value = 13812.5
$ws.Range("L68").Value = 2980
$ws.Range("M68").Value = -13063.5
$ws.Range("N68").Value = -4478
$ws.Range("H71").Value = 12608.889
$ws.Range("I71").Value = 13812.5
$ws.Range("J71").Value = 2980
$ws.Range("K71").Value = 69062.5
$ws.Range("L71").Value = 14900
$ws.Range("M71").Value = -65318.5
$ws.Range("N71").Value = -22388
$ws.Range("H82").Value = 1502.9688
$ws.Range("I82").Value = 1185.8889
$ws.Range("J82").Value = 1910.6428
$ws.Range("K82").Value = 1185.8889
$ws.Range("L82").Value = 1910.6428
$ws.Range("M82").Value = -824.8888999999999
$ws.Range("N82").Value = -2632.6428
$ws.Range("H85").Value = 1502.9688
$ws.Range("I85").Value = 1185.8889
$ws.Range("J85").Value = 1910.6428
$ws.Range("K85").Value = 1185.8889
$ws.Range("L85").Value = 1910.6428
$ws.Range("M85").Value = 62.11110000000008
$ws.Range("N85").Value = -4406.6428
$ws.Range("H94").Value = 72776.664
$ws.Range("J94").Value = 72776.664
$ws.Range("L94").Value = 72776.664
$ws.Range("N94").Value = -74128.664
$ws.Range("H136").Value = 264718.88
$ws.Range("J136").Value = 1937.2222
$ws.Range("L136").Value = 5811.6666
$ws.Range("N136").Value = -10911.6666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 45410
$ws.Range("J133").Value = 45410
$ws.Range("L133").Value = 45410
$ws.Range("N133").Value = -55530
